$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '30.907.85'
$ws.Cells.Item(2, 5).Value = '  -0.50%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.950.60'
$ws.Cells.Item(3, 5).Value = '  -0.70%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '242.09'
$ws.Cells.Item(5, 5).Value = '  -2.17%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.13%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4886'
$ws.Cells.Item(7, 5).Value = '  +0.37%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.82%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06965'
$ws.Cells.Item(9, 5).Value = '  +2.14%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '19.45'
$ws.Cells.Item(10, 5).Value = '  +1.26%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '106.89'
$ws.Cells.Item(11, 5).Value = '  -0.68%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.949.91'
$ws.Cells.Item(12, 5).Value = '  -0.73%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.07762'
$ws.Cells.Item(13, 5).Value = '  -0.15%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -1.64%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.6981'
$ws.Cells.Item(15, 5).Value = '  -0.60%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '278.97'
$ws.Cells.Item(16, 5).Value = '  -2.67%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '30.916.97'
$ws.Cells.Item(17, 5).Value = '  -0.45%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000007746'
$ws.Cells.Item(18, 5).Value = '  +0.12%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'Avalanche'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '13.17'
$ws.Cells.Item(19, 5).Value = '  -0.28%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '2.210.94'
$ws.Cells.Item(20, 5).Value = '  -0.62%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.002'
$ws.Cells.Item(21, 5).Value = '  +0.17%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.473'
$ws.Cells.Item(22, 5).Value = '  -2.38%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.003'
$ws.Cells.Item(23, 5).Value = '  +0.01%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '6.464'
$ws.Cells.Item(24, 5).Value = '  -2.09%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.718'
$ws.Cells.Item(25, 5).Value = '  -2.91%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '168.39'
$ws.Cells.Item(26, 5).Value = '  -0.73%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '19.64'
$ws.Cells.Item(27, 5).Value = '  -2.00%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '2.164'
$ws.Cells.Item(28, 5).Value = '  -1.10%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.1044'
$ws.Cells.Item(29, 5).Value = '  -1.93%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.398'
$ws.Cells.Item(30, 5).Value = '  -2.76%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '4.587'
$ws.Cells.Item(31, 5).Value = '  -5.31%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.555'
$ws.Cells.Item(32, 5).Value = '  -2.73%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.379'
$ws.Cells.Item(33, 5).Value = '  -2.78%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.04862'
$ws.Cells.Item(34, 5).Value = '  -4.44%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.7501'
$ws.Cells.Item(35, 5).Value = '  -2.66%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.162'
$ws.Cells.Item(36, 5).Value = '  -0.55%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.728'
$ws.Cells.Item(37, 5).Value = '  -0.07%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01992'
$ws.Cells.Item(38, 5).Value = '  -2.36%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.675'
$ws.Cells.Item(39, 5).Value = '  -1.95%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.494'
$ws.Cells.Item(40, 5).Value = '  +0.58%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '77.61'
$ws.Cells.Item(41, 5).Value = '  +6.84%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.093'
$ws.Cells.Item(42, 5).Value = '  -1.91%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.8936'
$ws.Cells.Item(43, 5).Value = '  +0.73%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '108.96'
$ws.Cells.Item(44, 5).Value = '  -0.86%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.4424'
$ws.Cells.Item(45, 5).Value = '  -1.16%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.9996'
$ws.Cells.Item(46, 5).Value = '  -0.04%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '7.778'
$ws.Cells.Item(47, 5).Value = '  +3.44%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '994.14'
$ws.Cells.Item(48, 5).Value = '  -0.75%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.1245'

# Row 50
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '9.215'
$ws.Cells.Item(50, 5).Value = '  -2.66%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Elrond'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '35.81'
$ws.Cells.Item(51, 5).Value = '  -0.47%  '
